$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 880.8329105377197
$ws.Range("E20").Value = 9.996771812438965
$ws.Range("E21").Value = 32.76991844177246
$ws.Range("E35").Value = 30.50720691680908
$ws.Range("E42").Value = 24.91283416748047
$ws.Range("E46").Value = 2.062201499938965
$ws.Range("E48").Value = 3.977298736572266
$ws.Range("E49").Value = 6.038784980773926
$ws.Range("E56").Value = 16.15941524505615
$ws.Range("E67").Value = 3.173708915710449
$ws.Range("E70").Value = 10.14125347137451
$ws.Range("E77").Value = 5.353331565856934
$ws.Range("E84").Value = 26.95858478546143
$ws.Range("E91").Value = 30.69961071014404
$ws.Range("E98").Value = 18.90277862548828
$ws.Range("E99").Value = 54.84485626220703
